$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 1.64
$ws.Cells.Item(2, 9).Value = 5.8
$ws.Cells.Item(2, 12).Value = 1.33
$ws.Cells.Item(2, 13).Value = 1.04
$ws.Cells.Item(2, 19).Value = 2.74
$ws.Cells.Item(2, 23).Value = 2.5
$ws.Cells.Item(2, 27).Value = 150

$ws.Cells.Item(3, 6).Value = 10
$ws.Cells.Item(3, 9).Value = 1.15
$ws.Cells.Item(3, 12).Value = 1.05
$ws.Cells.Item(3, 14).Value = 4.4
$ws.Cells.Item(3, 18).Value = 2.34
$ws.Cells.Item(3, 19).Value = 1.33
$ws.Cells.Item(3, 22).Value = 7.6

$ws.Cells.Item(4, 6).Value = 2.12
$ws.Cells.Item(4, 7).Value = 2.4
$ws.Cells.Item(4, 8).Value = 3.95
$ws.Cells.Item(4, 12).Value = 1.01

$ws.Cells.Item(5, 7).Value = 2.36
$ws.Cells.Item(5, 8).Value = 2.9
$ws.Cells.Item(5, 17).Value = 1.45
$ws.Cells.Item(5, 19).Value = 2.02
$ws.Cells.Item(5, 28).Value = 19
$ws.Cells.Item(5, 32).Value = 22
$ws.Cells.Item(5, 34).Value = 15
$ws.Cells.Item(5, 40).Value = 11

$ws.Cells.Item(6, 9).Value = 3.85
$ws.Cells.Item(6, 11).Value = 3.05
$ws.Cells.Item(6, 12).Value = 1.57
$ws.Cells.Item(6, 14).Value = 2.2
$ws.Cells.Item(6, 16).Value = 1.39
$ws.Cells.Item(6, 21).Value = 1.63
$ws.Cells.Item(6, 23).Value = 1.58
$ws.Cells.Item(6, 30).Value = 22

$ws.Cells.Item(7, 6).Value = 1.45
$ws.Cells.Item(7, 7).Value = 1.51
$ws.Cells.Item(7, 8).Value = 7.6
$ws.Cells.Item(7, 9).Value = 9.199999999999999
$ws.Cells.Item(7, 10).Value = 4.5
$ws.Cells.Item(7, 11).Value = 5.3
$ws.Cells.Item(7, 12).Value = 1.29
$ws.Cells.Item(7, 13).Value = 1.05
$ws.Cells.Item(7, 14).Value = 4.2
$ws.Cells.Item(7, 15).Value = 1.25
$ws.Cells.Item(7, 16).Value = 2.12
$ws.Cells.Item(7, 17).Value = 1.74
$ws.Cells.Item(7, 18).Value = 1.44
$ws.Cells.Item(7, 19).Value = 2.86
$ws.Cells.Item(7, 20).Value = 1.93
$ws.Cells.Item(7, 21).Value = 1.89
$ws.Cells.Item(7, 22).Value = 1.12
$ws.Cells.Item(7, 23).Value = 2.92
$ws.Cells.Item(7, 24).Value = 21
$ws.Cells.Item(7, 25).Value = 990
$ws.Cells.Item(7, 27).Value = 320
$ws.Cells.Item(7, 28).Value = 9
$ws.Cells.Item(7, 29).Value = 11.5
$ws.Cells.Item(7, 30).Value = 34
$ws.Cells.Item(7, 32).Value = 9.4
$ws.Cells.Item(7, 34).Value = 990
$ws.Cells.Item(7, 36).Value = 14
$ws.Cells.Item(7, 37).Value = 16.5
$ws.Cells.Item(7, 40).Value = 7.4

$ws.Cells.Item(8, 14).Value = 3.35
$ws.Cells.Item(8, 17).Value = 2.02
$ws.Cells.Item(8, 18).Value = 1.31
$ws.Cells.Item(8, 19).Value = 3.6
$ws.Cells.Item(8, 20).Value = 1.77
$ws.Cells.Item(8, 21).Value = 2.06
$ws.Cells.Item(8, 24).Value = 16
$ws.Cells.Item(8, 25).Value = 13
$ws.Cells.Item(8, 32).Value = 24
$ws.Cells.Item(8, 34).Value = 22

$ws.Cells.Item(9, 7).Value = 1.63
$ws.Cells.Item(9, 8).Value = 8.6
$ws.Cells.Item(9, 11).Value = 3.7
$ws.Cells.Item(9, 15).Value = 1.64
$ws.Cells.Item(9, 23).Value = 2.58

$ws.Cells.Item(10, 6).Value = 1.7
$ws.Cells.Item(10, 7).Value = 1.71
$ws.Cells.Item(10, 9).Value = 7
$ws.Cells.Item(10, 12).Value = 1.55
$ws.Cells.Item(10, 15).Value = 1.49
$ws.Cells.Item(10, 16).Value = 1.65
$ws.Cells.Item(10, 17).Value = 2.48
$ws.Cells.Item(10, 18).Value = 1.23
$ws.Cells.Item(10, 20).Value = 2.34
$ws.Cells.Item(10, 22).Value = 1.16
$ws.Cells.Item(10, 23).Value = 2.4
$ws.Cells.Item(10, 27).Value = 230
$ws.Cells.Item(10, 28).Value = 6.2
$ws.Cells.Item(10, 29).Value = 8.4
$ws.Cells.Item(10, 30).Value = 27
$ws.Cells.Item(10, 35).Value = 150
$ws.Cells.Item(10, 39).Value = 230

$ws.Cells.Item(11, 6).Value = 1.11
$ws.Cells.Item(11, 9).Value = 42
$ws.Cells.Item(11, 11).Value = 12.5
$ws.Cells.Item(11, 16).Value = 2.9
$ws.Cells.Item(11, 17).Value = 1.47
$ws.Cells.Item(11, 18).Value = 1.73
$ws.Cells.Item(11, 19).Value = 2.22
$ws.Cells.Item(11, 21).Value = 1.53
$ws.Cells.Item(11, 24).Value = 46
$ws.Cells.Item(11, 29).Value = 990
